# Applies the diff: turns the two formerly-bold content paragraphs
# ("-Función innecesaria." / "-Optimización ...") into bold section
# headers ("Diagrama:" / "Caja Negra:"), moving their old (now
# non-bold) text into freshly inserted paragraphs right below the new
# "Código:" heading, and appends several brand-new plain-text
# (non-bold) paragraphs of commentary. The "_GoBack" bookmark is
# re-created at the end of the document's final paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Paragraph 2 was empty (bold pPr) -> becomes "Código:" (bold).
# ---------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter("Código:")
$p2.Range.Font.Bold = $True

# ---------------------------------------------------------------
# 2) Insert three new (non-bold) paragraphs right after "Código:" and
#    before the old "-Función innecesaria." paragraph:
#      - "-Función innecesaria."
#      - "-Optimización "admitirpalabra"."
#      - "-No se pide la palabra a la hora de jugar. Solo la posición
#        y la orientación."
# ---------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$pFuncion = $d.Paragraphs(3)
$pFuncion.Range.InsertAfter("-Función innecesaria.")
$pFuncion.Range.Font.Bold = $False

$pFuncion.Range.InsertParagraphAfter()
$pOptim = $d.Paragraphs(4)
$pOptim.Range.InsertAfter("-Optimización")
$pOptim.Range.InsertAfter(" “")
$pOptim.Range.InsertAfter("admitirpalabra")
$pOptim.Range.InsertAfter("”.")
$pOptim.Range.Font.Bold = $False

$pOptim.Range.InsertParagraphAfter()
$pNoPide = $d.Paragraphs(5)
$pNoPide.Range.InsertAfter("-No se pide la palabra a la hora de jugar.")
$pNoPide.Range.InsertAfter(" Solo la posición y la orientación.")
$pNoPide.Range.Font.Bold = $False

# ---------------------------------------------------------------
# 3) The old "-Función innecesaria." paragraph (now paragraph 6) keeps
#    its bold paragraph formatting but its text becomes "Diagrama:".
# ---------------------------------------------------------------
$pDiagrama = $d.Paragraphs(6)
$r = $pDiagrama.Range
$r.End = $r.End - 1
$r.Text = "Diagrama:"
$pDiagrama.Range.Font.Bold = $True

# ---------------------------------------------------------------
# 4) Insert two new (non-bold) paragraphs after "Diagrama:", plus one
#    empty paragraph:
#      - "Dibujar sopa de letras incluye crear sopa de letras."
#      - "Jugar sopa de letras incluye crear sopa de letras."
#      - "" (empty)
# ---------------------------------------------------------------
$pDiagrama.Range.InsertParagraphAfter()
$pDibujar = $d.Paragraphs(7)
$pDibujar.Range.InsertAfter("Dibujar")
$pDibujar.Range.InsertAfter(" sopa de letras incluye crear sopa de letras")
$pDibujar.Range.InsertAfter(".")
$pDibujar.Range.Font.Bold = $False

$pDibujar.Range.InsertParagraphAfter()
$pJugar = $d.Paragraphs(8)
$pJugar.Range.InsertAfter("Jugar sopa de letras incluye crear sopa de letras.")
$pJugar.Range.Font.Bold = $False

$pJugar.Range.InsertParagraphAfter()

# ---------------------------------------------------------------
# 5) The old "-Optimización ..." paragraph (now paragraph 10) keeps
#    its bold paragraph formatting but its text becomes "Caja Negra:".
#    The trailing text and the mid-paragraph bookmark that used to
#    follow "-Optimización" are dropped here; the bookmark is
#    re-created at the end of the new final paragraph below.
# ---------------------------------------------------------------
$pCaja = $d.Paragraphs(10)
$r = $pCaja.Range
$r.End = $r.End - 1
$r.Text = "Caja Negra:"
$pCaja.Range.Font.Bold = $True

# ---------------------------------------------------------------
# 6) Insert the final (non-bold) paragraph with the "admitir palabra"
#    commentary, and re-create the _GoBack bookmark at its end.
# ---------------------------------------------------------------
$pCaja.Range.InsertParagraphAfter()
$pAdmitir = $d.Paragraphs(11)
$pAdmitir.Range.InsertAfter("“")
$pAdmitir.Range.InsertAfter("a")
$pAdmitir.Range.InsertAfter("dmitir")
$pAdmitir.Range.InsertAfter("pala")
$pAdmitir.Range.InsertAfter("bra")
$pAdmitir.Range.InsertAfter("”")
$pAdmitir.Range.InsertAfter(": admitir palabra que no debe. (no válido)")
$pAdmitir.Range.Font.Bold = $False

$bmRange = $pAdmitir.Range
$bmRange.End = $bmRange.End - 1
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
